$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "yuva"
$ws.Range("B3").Value = "yuva"
$ws.Range("C3").Value = "Automation Anywhere"

$ws.Range("A4").Select()
